# Append two new listings, insert them above the two existing bottom rows,
# and append two more at the end; refresh the "fetched at" timestamp on
# every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2026-02-06 01:52:51"

# --- refresh timestamp on all existing data rows (2-10) ---
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = $timestamp
}

# --- make room: push the current rows 9-10 down to 11-12 ---
$ws.Rows.Item(11).Resize(2, 1).EntireRow.Insert()

$ws.Cells.Item(11, 1).Value = $timestamp
$ws.Cells.Item(11, 2).Value = "初回 line予約システム、Googlrカレンダー連動一元管理"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5486342"
$ws.Cells.Item(11, 7).Value = 53
$ws.Cells.Item(11, 8).Value = "◇管理"

$ws.Cells.Item(12, 1).Value = $timestamp
$ws.Cells.Item(12, 2).Value = "【急募】iOS/AndroidアプリのSkyWay切替対応エンジニア募集"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5486110"
$ws.Cells.Item(12, 7).Value = 38
$ws.Cells.Item(12, 8).Value = "◇アプリ"

# --- overwrite row 9 with the new top listing ---
$ws.Cells.Item(9, 1).Value = $timestamp
$ws.Cells.Item(9, 2).Value = "【業務改善】訪問業務に特化したスケジュール/介入実績管理Webシステム構築"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5486583"
$ws.Cells.Item(9, 7).Value = 85
$ws.Cells.Item(9, 8).Value = "◇業務改善"

# --- overwrite row 10 with another new listing ---
$ws.Cells.Item(10, 1).Value = $timestamp
$ws.Cells.Item(10, 2).Value = "【Java/講師】新入社員研修のサブ講師募集"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5477144"
$ws.Cells.Item(10, 7).Value = 78
$ws.Cells.Item(10, 8).Value = "★Java"

# --- append two brand-new rows at the bottom ---
$ws.Cells.Item(13, 1).Value = $timestamp
$ws.Cells.Item(13, 2).Value = "【長期】寝具ブランドのAmazon・楽天市場 運用代行パートナー募集"
$ws.Cells.Item(13, 3).Value = "システム開発"
$ws.Cells.Item(13, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(13, 5).Value = "期限情報なし"
$ws.Cells.Item(13, 6).Value = "https://www.lancers.jp/work/detail/5486471"
$ws.Cells.Item(13, 7).Value = 25

$ws.Cells.Item(14, 1).Value = $timestamp
$ws.Cells.Item(14, 2).Value = "【急募】Klaviyoスパム対策とドメイン解決の専門家募集"
$ws.Cells.Item(14, 3).Value = "システム開発"
$ws.Cells.Item(14, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(14, 5).Value = "期限情報なし"
$ws.Cells.Item(14, 6).Value = "https://www.lancers.jp/work/detail/5486673"
$ws.Cells.Item(14, 7).Value = 13

# --- hyperlinks for the URL cells that are new/changed (rows 9-14); the
#     existing rows 2-8 already carry a working hyperlink and are left as-is ---
for ($r = 9; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value
    $ws.Hyperlinks.Add($cell, $url) | Out-Null
    $cell.Style = "Hyperlink"
}
